# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" and refresh the
# related handoff timestamps across the Overview, zh-cn and de-de sheets.
# Also widen the "Status" column on each sheet to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Sheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-17 08:38:41"
# ColumnWidth is quantized to pixel/character increments by Excel, so 16.33
# is the nearest character-width setting that lands on the saved column
# width closest to the target (~17.22 "file units").
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-17 08:38:36"
$zhcn.Columns.Item(3).ColumnWidth = 16.33

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Sheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-17 08:38:41"
$dede.Columns.Item(3).ColumnWidth = 16.33
